# Fruta / hortaliza, semanal
# Insert a new weekly record row right after existing row 330 (i.e. as the new
# row 331), shifting all the subsequent rows (old 331..384) down by one to
# 332..385, and populate the new row with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 331; existing row 331 (and below) shift to 332+.
$ws.Rows.Item(331).Insert()

# Populate the newly inserted row 331 with the new weekly record.
$ws.Range("A331").Value = 6
$ws.Range("B331").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C331").Value = "Metropolitana"
$ws.Range("D331").Value = "7/28/2023"
$ws.Range("E331").Value = 13
$ws.Range("F331").Value = "Fruta"
$ws.Range("G331").Value = 100107
$ws.Range("H331").Value = "Otros"
$ws.Range("I331").Value = 100107002
$ws.Range("J331").Value = "Chirimoya"
$ws.Range("K331").Value = "Cultivar IV Región"
$ws.Range("L331").Value = "Segunda"
$ws.Range("M331").Value = 200
$ws.Range("N331").Value = 26000
$ws.Range("O331").Value = 26000
$ws.Range("P331").Value = 26000
$ws.Range("Q331").Value = "$/bandeja 10 kilos"
$ws.Range("R331").Value = "Provincia de Limarí"
$ws.Range("S331").Value = 2600
$ws.Range("T331").Value = 10
